# The workbook's single sheet is a flat data table (header row 1, data
# rows 2-494). The edit inserts one new data record for "Vega Modelo de
# Temuco" / "Coliflor" at row 382, pushing the former rows 382-494 down by
# one (to 383-495) and growing the sheet dimension from R494 to R495.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 382, shifting rows 382..494 down to
# 383..495 (same as Excel's Home > Insert > Insert Sheet Rows on row 382).
$ws.Rows.Item(382).Insert()

# Populate the newly-inserted row 382 with the new record's values.
$ws.Cells.Item(382, 1).Value  = 10
$ws.Cells.Item(382, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(382, 3).Value  = "La Araucanía"
$ws.Cells.Item(382, 4).Value  = 44876
$ws.Cells.Item(382, 5).Value  = 9
$ws.Cells.Item(382, 6).Value  = 100112008
$ws.Cells.Item(382, 7).Value  = "Coliflor"
$ws.Cells.Item(382, 8).Value  = "Sin especificar"
$ws.Cells.Item(382, 9).Value  = "Primera"
$ws.Cells.Item(382, 10).Value = 2150
$ws.Cells.Item(382, 11).Value = 900
$ws.Cells.Item(382, 12).Value = 1000
$ws.Cells.Item(382, 13).Value = 942
$ws.Cells.Item(382, 14).Value = "$/unidad"
$ws.Cells.Item(382, 15).Value = "Región del Maule"
$ws.Cells.Item(382, 16).Value = 942
$ws.Cells.Item(382, 17).Value = 1
$ws.Cells.Item(382, 18).Value = "Hortaliza"
